$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 5: split the former "Java와 Mysql 연동하기" task into a numbered
# sub-task and shift the result/problem columns over.
$ws.Range("A5").Value = "Java와 Mysql 연동하기(1)"
$ws.Range("E5").Value = "testing 코드 git commit"
$ws.Range("F5").Value = "실제 gui java 파일에서 tomcat server가 연결문제 처리하는데 시간 필요"

# Row 6: new task entry
$ws.Range("A6").Value = "Java와 Mysql 연동하기(2)"
$ws.Range("C6").Value = "2019-05-17"
$ws.Range("D6").Value = "2019-05-18"
$ws.Range("E6").Value = "수정한 Storage.java git commit"

# Row 7: new task entry
$ws.Range("A7").Value = "Java와 Mysql 연동하기(3)"
$ws.Range("B7").Value = "Member.java와 MySql 연동하기"
$ws.Range("C7").Value = "2019-05-18"

# Row 6 content column (written after row 7 to reproduce the author's
# original shared-string ordering)
$ws.Range("B6").Value = "Storage.java와 Mysql 연동하기"

# Row 8: new task entry
$ws.Range("A8").Value = "Java와 Mysql 연동하기(4)"
$ws.Range("B8").Value = "Tabel.java, Menu.java MySql 연동하기"
$ws.Range("C8").Value = "2019-05-18"

# Row heights to match the wrapped text content
$ws.Rows.Item(6).RowHeight = 54.85
$ws.Rows.Item(7).RowHeight = 57.4
$ws.Rows.Item(8).RowHeight = 76.5

# Update the selection / active cell shown in the sheet view
$ws.Range("H8").Select()
